$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "N" column formulas so the Soli tax is expressed in monthly
# terms (divide the existing yearly computation by 12).
$ws.Range("N2").Formula = "=MIN(0.055*P2,MAX(0.2*P2,972))/12"
$ws.Range("N3").Formula = "=MIN(0.055*P3,MAX(0.2*P3,972))/12"
$ws.Range("N5").Formula = "=MIN(0.055*P5,MAX(0.2*P5,972))/12"
$ws.Range("N6").Formula = "=MIN(0.055*P6,MAX(0.2*P6,972))/12"
$ws.Range("N7").Formula = "=MIN(0.055*P7,MAX(0.2*P7,972))/12"

# Move the active selection to N6, matching the saved view state.
$ws.Range("N6").Select()
